# Update Betfair Back/Lay odds for 2026-01-02 (row 2, 5, 6, 7, 8)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Australian A-League Men: Melbourne Victory vs Perth Glory
$ws.Range("H2").Value = 4.4
$ws.Range("W2").Value = 2.12
$ws.Range("AG2").Value = 9.6

# Row 5 - Saudi Professional League: Al Ahli vs Al Nassr
$ws.Range("F5").Value = 3.8
$ws.Range("G5").Value = 4.5
$ws.Range("H5").Value = 1.79
$ws.Range("I5").Value = 1.96
$ws.Range("J5").Value = 4.1
$ws.Range("K5").Value = 4.9
$ws.Range("Q5").Value = 1.52

# Row 6 - French Ligue 1: Toulouse vs Lens
$ws.Range("F6").Value = 2.9
$ws.Range("I6").Value = 2.72
$ws.Range("T6").Value = 1.77
$ws.Range("U6").Value = 2.2

# Row 7 - Italian Serie A: Cagliari vs AC Milan
$ws.Range("P7").Value = 1.99
$ws.Range("T7").Value = 2.06

# Row 8 - Spanish La Liga: Rayo Vallecano vs Getafe
$ws.Range("T8").Value = 2.66
$ws.Range("AD8").Value = 21
$ws.Range("AJ8").Value = 38
$ws.Range("AK8").Value = 42
